$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 1398.5
$ws.Range("I18").Value = 1398.5
$ws.Range("K18").Value = 1398.5
$ws.Range("M18").Value = -1114.5
$ws.Range("H38").Value = 590.3333
$ws.Range("I38").Value = 590.3333
$ws.Range("K38").Value = 1770.9999
$ws.Range("M38").Value = -1398.9999
$ws.Range("H40").Value = 3374.75
$ws.Range("I40").Value = 0
$ws.Range("K40").Value = 0
$ws.Range("M40").ClearContents()
$ws.Range("H43").Value = 9999.166999999999
$ws.Range("I43").Value = 9998.333000000001
$ws.Range("J43").Value = 10000
$ws.Range("K43").Value = 9998.333000000001
$ws.Range("L43").Value = 10000
$ws.Range("M43").Value = -9929.333000000001
$ws.Range("N43").Value = -10138
$ws.Range("H58").Value = 1979.9333
$ws.Range("J58").Value = 10017
$ws.Range("L58").Value = 30051
$ws.Range("N58").Value = -30351
$ws.Range("H61").Value = 297.8
$ws.Range("I61").Value = 297.8
$ws.Range("K61").Value = 893.4000000000001
$ws.Range("M61").Value = -721.4000000000001
$ws.Range("H62").Value = 10085.19
$ws.Range("I62").Value = 9873.916999999999
$ws.Range("J62").Value = 10366.889
$ws.Range("K62").Value = 9873.916999999999
$ws.Range("L62").Value = 10366.889
$ws.Range("M62").Value = -9249.916999999999
$ws.Range("N62").Value = -11614.889
$ws.Range("H65").Value = 10085.19
$ws.Range("I65").Value = 9873.916999999999
$ws.Range("J65").Value = 10366.889
$ws.Range("K65").Value = 49369.585
$ws.Range("L65").Value = 51834.44499999999
$ws.Range("M65").Value = -46249.585
$ws.Range("N65").Value = -58074.44499999999
$ws.Range("H69").Value = 459299.6
$ws.Range("I69").Value = 9000
$ws.Range("J69").Value = 509332.88
$ws.Range("K69").Value = 27000
$ws.Range("L69").Value = 1527998.64
$ws.Range("M69").Value = -26126
$ws.Range("N69").Value = -1529746.64
$ws.Range("H70").Value = 38296.78
$ws.Range("I70").Value = 6288
$ws.Range("J70").Value = 52846.227
$ws.Range("K70").Value = 18864
$ws.Range("L70").Value = 158538.681
$ws.Range("M70").Value = -18594
$ws.Range("N70").Value = -159078.681
$ws.Range("H72").Value = 459299.6
$ws.Range("I72").Value = 9000
$ws.Range("J72").Value = 509332.88
$ws.Range("K72").Value = 81000
$ws.Range("L72").Value = 4583995.92
$ws.Range("M72").Value = -76632
$ws.Range("N72").Value = -4592731.92
$ws.Range("H73").Value = 38296.78
$ws.Range("I73").Value = 6288
$ws.Range("J73").Value = 52846.227
$ws.Range("K73").Value = 18864
$ws.Range("L73").Value = 158538.681
$ws.Range("M73").Value = -17928
$ws.Range("N73").Value = -160410.681
$ws.Range("H76").Value = 6805.278
$ws.Range("I76").Value = 6292.75
$ws.Range("J76").Value = 7215.3
$ws.Range("K76").Value = 6292.75
$ws.Range("L76").Value = 7215.3
$ws.Range("M76").Value = -5977.75
$ws.Range("N76").Value = -7845.3
$ws.Range("H79").Value = 6805.278
$ws.Range("I79").Value = 6292.75
$ws.Range("J79").Value = 7215.3
$ws.Range("K79").Value = 6292.75
$ws.Range("L79").Value = 7215.3
$ws.Range("M79").Value = -5200.75
$ws.Range("N79").Value = -9399.299999999999
$ws.Range("H80").Value = 147902.17
$ws.Range("I80").Value = 8261.076999999999
$ws.Range("J80").Value = 312932.53
$ws.Range("K80").Value = 24783.231
$ws.Range("L80").Value = 938797.5900000001
$ws.Range("M80").Value = -23785.231
$ws.Range("N80").Value = -940793.5900000001
$ws.Range("H82").Value = 1478.091
$ws.Range("I82").Value = 625.9
$ws.Range("K82").Value = 1877.7
$ws.Range("M82").Value = -1471.7
$ws.Range("H83").Value = 147902.17
$ws.Range("I83").Value = 8261.076999999999
$ws.Range("J83").Value = 312932.53
$ws.Range("K83").Value = 74349.693
$ws.Range("L83").Value = 2816392.77
$ws.Range("M83").Value = -69357.693
$ws.Range("N83").Value = -2826376.77
$ws.Range("H85").Value = 1478.091
$ws.Range("I85").Value = 625.9
$ws.Range("K85").Value = 1877.7
$ws.Range("M85").Value = -473.6999999999998
$ws.Range("H86").Value = 5521.8
$ws.Range("I86").Value = 4705
$ws.Range("J86").Value = 6066.3335
$ws.Range("K86").Value = 4705
$ws.Range("L86").Value = 6066.3335
$ws.Range("M86").Value = -3582
$ws.Range("N86").Value = -8312.333500000001
$ws.Range("H89").Value = 5521.8
$ws.Range("I89").Value = 4705
$ws.Range("J89").Value = 6066.3335
$ws.Range("K89").Value = 23525
$ws.Range("L89").Value = 30331.6675
$ws.Range("M89").Value = -17909
$ws.Range("N89").Value = -41563.6675
$ws.Range("H92").Value = 2381.1428
$ws.Range("I92").Value = 423.6
$ws.Range("J92").Value = 7275
$ws.Range("K92").Value = 423.6
$ws.Range("L92").Value = 7275
$ws.Range("M92").Value = 824.4
$ws.Range("N92").Value = -9771
$ws.Range("H96").Value = 8929418
$ws.Range("I96").Value = 17857290
$ws.Range("J96").Value = 1547.25
$ws.Range("K96").Value = 53571870
$ws.Range("L96").Value = 4641.75
$ws.Range("M96").Value = -53570497
$ws.Range("N96").Value = -7387.75
$ws.Range("H112").Value = 8159.024
$ws.Range("J112").Value = 8653.59
$ws.Range("L112").Value = 25960.77
$ws.Range("N112").Value = -28176.77
$ws.Range("H116").Value = 6570.2383
$ws.Range("J116").Value = 4799.2
$ws.Range("L116").Value = 4799.2
$ws.Range("N116").Value = -11683.2
$ws.Range("H125").Value = 23216.572
$ws.Range("I125").Value = 38699.875
$ws.Range("J125").Value = 2572.1667
$ws.Range("K125").Value = 348298.875
$ws.Range("L125").Value = 23149.5003
$ws.Range("M125").Value = -345838.875
$ws.Range("N125").Value = -28069.5003
$ws.Range("H127").Value = 1116.7391
$ws.Range("I127").Value = 809.5789
$ws.Range("J127").Value = 2575.75
$ws.Range("K127").Value = 2428.7367
$ws.Range("L127").Value = 7727.25
$ws.Range("M127").Value = 2531.2633
$ws.Range("N127").Value = -17647.25
$ws.Range("H132").Value = 1994714.8
$ws.Range("I132").Value = 1994714.8
$ws.Range("K132").Value = 5984144.4
$ws.Range("M132").Value = -5981614.4
$ws.Range("H133").Value = 150000
$ws.Range("J133").Value = 150000
$ws.Range("L133").Value = 150000
$ws.Range("N133").Value = -160120
$ws.Range("H135").Value = 14388.841
$ws.Range("I135").Value = 1119.8572
$ws.Range("K135").Value = 10078.7148
$ws.Range("M135").Value = -7543.7148
$ws.Range("H137").Value = 15358.444
$ws.Range("J137").Value = 4344.154
$ws.Range("L137").Value = 13032.462
$ws.Range("N137").Value = -18132.462
$ws.Range("H138").Value = 3686.2654
$ws.Range("J138").Value = 3850.8857
$ws.Range("L138").Value = 11552.6571
$ws.Range("N138").Value = -21832.6571
$ws.Range("H141").Value = 1306.8846
$ws.Range("I141").Value = 1306.8846
$ws.Range("K141").Value = 3920.6538
$ws.Range("M141").Value = 1259.3462

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 1049.6666
$ws.Range("I4").Value = 1499
$ws.Range("K4").Value = 1499
$ws.Range("M4").Value = -1383
$ws.Range("H32").Value = 8253.823
$ws.Range("I32").Value = 8083.936
$ws.Range("J32").Value = 10250
$ws.Range("K32").Value = 8083.936
$ws.Range("L32").Value = 10250
$ws.Range("M32").Value = -7796.936
$ws.Range("N32").Value = -10824
$ws.Range("H38").Value = 10508.5
$ws.Range("J38").Value = 0
$ws.Range("L38").Value = 0
$ws.Range("N38").ClearContents()
$ws.Range("H61").Value = 2819.6558
$ws.Range("I61").Value = 2144.375
$ws.Range("J61").Value = 10382.8
$ws.Range("K61").Value = 2144.375
$ws.Range("L61").Value = 10382.8
$ws.Range("M61").Value = -1932.375
$ws.Range("N61").Value = -10806.8
$ws.Range("H63").Value = 3410.5
$ws.Range("I63").Value = 3410.5
$ws.Range("K63").Value = 3410.5
$ws.Range("M63").Value = -2724.5
$ws.Range("H66").Value = 3410.5
$ws.Range("I66").Value = 3410.5
$ws.Range("K66").Value = 17052.5
$ws.Range("M66").Value = -13620.5
$ws.Range("H74").Value = 1004333.8
$ws.Range("I74").Value = 2002668.6
$ws.Range("K74").Value = 2002668.6
$ws.Range("M74").Value = -2001794.6
$ws.Range("H76").Value = 400000
$ws.Range("J76").Value = 400000
$ws.Range("L76").Value = 400000
$ws.Range("N76").Value = -400676
$ws.Range("H77").Value = 1004333.8
$ws.Range("I77").Value = 2002668.6
$ws.Range("K77").Value = 10013343
$ws.Range("M77").Value = -10008975
$ws.Range("H79").Value = 400000
$ws.Range("J79").Value = 400000
$ws.Range("L79").Value = 400000
$ws.Range("N79").Value = -402340
$ws.Range("H88").Value = 5181.1333
$ws.Range("I88").Value = 1084.5
$ws.Range("K88").Value = 1084.5
$ws.Range("M88").Value = -678.5
$ws.Range("H91").Value = 5181.1333
$ws.Range("I91").Value = 1084.5
$ws.Range("K91").Value = 1084.5
$ws.Range("M91").Value = 319.5
$ws.Range("H97").Value = 3356.6667
$ws.Range("I97").Value = 1713.3334
$ws.Range("K97").Value = 1713.3334
$ws.Range("M97").Value = -1217.3334
$ws.Range("H119").Value = 148698
$ws.Range("J119").Value = 148698
$ws.Range("L119").Value = 148698
$ws.Range("N119").Value = -158374
$ws.Range("H122").Value = 2931.5417
$ws.Range("I122").Value = 2830.1667
$ws.Range("K122").Value = 8490.500100000001
$ws.Range("M122").Value = -6040.500100000001
$ws.Range("H126").Value = 5624.25
$ws.Range("I126").Value = 5624.25
$ws.Range("K126").Value = 16872.75
$ws.Range("M126").Value = -14402.75
$ws.Range("H132").Value = 1346.2122
$ws.Range("I132").Value = 1299.4746
$ws.Range("J132").Value = 1740.1428
$ws.Range("K132").Value = 3898.4238
$ws.Range("L132").Value = 5220.428400000001
$ws.Range("M132").Value = -1368.4238
$ws.Range("N132").Value = -10280.4284
$ws.Range("H136").Value = 2819.6558
$ws.Range("I136").Value = 2144.375
$ws.Range("J136").Value = 10382.8
$ws.Range("K136").Value = 6433.125
$ws.Range("L136").Value = 31148.4
$ws.Range("M136").Value = -3883.125
$ws.Range("N136").Value = -36248.39999999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1988.8572
$ws.Range("I3").Value = 1765.5927
$ws.Range("K3").Value = 1765.5927
$ws.Range("M3").Value = -1651.5927
$ws.Range("H20").Value = 301354.4
$ws.Range("I20").Value = 251691
$ws.Range("K20").Value = 251691
$ws.Range("M20").Value = -251444
$ws.Range("H22").Value = 708.5294
$ws.Range("I22").Value = 395.84616
$ws.Range("J22").Value = 1724.75
$ws.Range("K22").Value = 395.84616
$ws.Range("L22").Value = 1724.75
$ws.Range("M22").Value = -222.84616
$ws.Range("N22").Value = -2070.75
$ws.Range("H86").Value = 2000
$ws.Range("I86").Value = 2000
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 2000
$ws.Range("L86").Value = 0
$ws.Range("M86").Value = -877
$ws.Range("N86").ClearContents()
$ws.Range("H89").Value = 2000
$ws.Range("I89").Value = 2000
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 10000
$ws.Range("L89").Value = 0
$ws.Range("M89").Value = -4384
$ws.Range("N89").ClearContents()
$ws.Range("H94").Value = 1569.6666
$ws.Range("I94").Value = 1354.5
$ws.Range("K94").Value = 1354.5
$ws.Range("M94").Value = -903.5
$ws.Range("H99").Value = 7275.5
$ws.Range("I99").Value = 8360
$ws.Range("K99").Value = 8360
$ws.Range("M99").Value = -6862
$ws.Range("H105").Value = 3849.75
$ws.Range("I105").Value = 3749.5
$ws.Range("K105").Value = 3749.5
$ws.Range("M105").Value = -2002.5
$ws.Range("H128").Value = 2099
$ws.Range("I128").Value = 2099
$ws.Range("K128").Value = 6297
$ws.Range("M128").Value = -3807
$ws.Range("H134").Value = 2846.75
$ws.Range("I134").Value = 2601.3547
$ws.Range("K134").Value = 7804.0641
$ws.Range("M134").Value = -5269.0641

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 281.92856
$ws.Range("I7").Value = 285.5
$ws.Range("J7").Value = 273
$ws.Range("K7").Value = 285.5
$ws.Range("L7").Value = 273
$ws.Range("M7").Value = -172.5
$ws.Range("N7").Value = -499
$ws.Range("H28").Value = 11997.833
$ws.Range("J28").Value = 11997.833
$ws.Range("L28").Value = 11997.833
$ws.Range("N28").Value = -12487.833
$ws.Range("H31").Value = 5949.56
$ws.Range("I31").Value = 7504
$ws.Range("J31").Value = 5345.0557
$ws.Range("K31").Value = 7504
$ws.Range("L31").Value = 5345.0557
$ws.Range("M31").Value = -7209
$ws.Range("N31").Value = -5935.0557
$ws.Range("H34").Value = 5949.56
$ws.Range("I34").Value = 7504
$ws.Range("J34").Value = 5345.0557
$ws.Range("K34").Value = 7504
$ws.Range("L34").Value = 5345.0557
$ws.Range("M34").Value = -7302
$ws.Range("N34").Value = -5749.0557
$ws.Range("H41").Value = 9909
$ws.Range("J41").Value = 30999.666
$ws.Range("L41").Value = 30999.666
$ws.Range("N41").Value = -31855.666
$ws.Range("H47").Value = 2249.5
$ws.Range("J47").Value = 2749.3333
$ws.Range("L47").Value = 2749.3333
$ws.Range("N47").Value = -3881.3333
$ws.Range("H48").Value = 9999
$ws.Range("I48").Value = 9999
$ws.Range("K48").Value = 9999
$ws.Range("M48").Value = -9523
$ws.Range("H51").Value = 44616.668
$ws.Range("J51").Value = 44616.668
$ws.Range("L51").Value = 44616.668
$ws.Range("N51").Value = -46088.668
$ws.Range("H58").Value = 3160.25
$ws.Range("I58").Value = 3281.4736
$ws.Range("J58").Value = 2831.2144
$ws.Range("K58").Value = 3281.4736
$ws.Range("L58").Value = 2831.2144
$ws.Range("M58").Value = -3078.4736
$ws.Range("N58").Value = -3237.2144
$ws.Range("H60").Value = 49249.93
$ws.Range("I60").Value = 29500
$ws.Range("K60").Value = 29500
$ws.Range("M60").Value = -28989
$ws.Range("H61").Value = 44616.668
$ws.Range("J61").Value = 44616.668
$ws.Range("L61").Value = 44616.668
$ws.Range("N61").Value = -45312.668
$ws.Range("H62").Value = 22011.111
$ws.Range("J62").Value = 30500
$ws.Range("L62").Value = 30500
$ws.Range("N62").Value = -31748
$ws.Range("H65").Value = 22011.111
$ws.Range("J65").Value = 30500
$ws.Range("L65").Value = 152500
$ws.Range("N65").Value = -158740
$ws.Range("H74").Value = 23499.5
$ws.Range("J74").Value = 0
$ws.Range("L74").Value = 0
$ws.Range("N74").ClearContents()
$ws.Range("H77").Value = 23499.5
$ws.Range("J77").Value = 0
$ws.Range("L77").Value = 0
$ws.Range("N77").ClearContents()
$ws.Range("H80").Value = 34999.8
$ws.Range("J80").Value = 34999.8
$ws.Range("L80").Value = 34999.8
$ws.Range("N80").Value = -37245.8
$ws.Range("H83").Value = 34999.8
$ws.Range("J83").Value = 34999.8
$ws.Range("L83").Value = 104999.4
$ws.Range("N83").Value = -116231.4
$ws.Range("H88").Value = 44999.5
$ws.Range("J88").Value = 44999.5
$ws.Range("L88").Value = 44999.5
$ws.Range("N88").Value = -45811.5
$ws.Range("H91").Value = 44999.5
$ws.Range("J91").Value = 44999.5
$ws.Range("L91").Value = 44999.5
$ws.Range("N91").Value = -47807.5
$ws.Range("H94").Value = 2644.5715
$ws.Range("I94").Value = 2778.25
$ws.Range("J94").Value = 2466.3333
$ws.Range("K94").Value = 2778.25
$ws.Range("L94").Value = 2466.3333
$ws.Range("M94").Value = -2327.25
$ws.Range("N94").Value = -3368.3333
$ws.Range("H95").Value = 18750
$ws.Range("J95").Value = 18750
$ws.Range("L95").Value = 18750
$ws.Range("N95").Value = -24242
$ws.Range("H107").Value = 1480.8182
$ws.Range("I107").Value = 535.9286
$ws.Range("J107").Value = 3134.375
$ws.Range("K107").Value = 535.9286
$ws.Range("L107").Value = 3134.375
$ws.Range("M107").Value = 1384.0714
$ws.Range("N107").Value = -6974.375
$ws.Range("H122").Value = 5958.353
$ws.Range("I122").Value = 5614.3105
$ws.Range("K122").Value = 16842.9315
$ws.Range("M122").Value = -14392.9315
$ws.Range("H132").Value = 31724.7
$ws.Range("I132").Value = 33341.74
$ws.Range("K132").Value = 100025.22
$ws.Range("M132").Value = -97495.22
$ws.Range("H133").Value = 159919.4
$ws.Range("J133").Value = 159825.25
$ws.Range("L133").Value = 159825.25
$ws.Range("N133").Value = -164885.25
$ws.Range("H134").Value = 3175.0952
$ws.Range("I134").Value = 2046.2
$ws.Range("K134").Value = 6138.6
$ws.Range("M134").Value = -3603.6
$ws.Range("H136").Value = 3160.25
$ws.Range("I136").Value = 3281.4736
$ws.Range("J136").Value = 2831.2144
$ws.Range("K136").Value = 9844.4208
$ws.Range("L136").Value = 8493.643199999999
$ws.Range("M136").Value = -7294.4208
$ws.Range("N136").Value = -13593.6432
$ws.Range("H137").Value = 91995
$ws.Range("J137").Value = 91995
$ws.Range("L137").Value = 91995
$ws.Range("N137").Value = -102195
$ws.Range("H139").Value = 61597.8
$ws.Range("I139").Value = 50000
$ws.Range("J139").Value = 69329.664
$ws.Range("K139").Value = 50000
$ws.Range("L139").Value = 69329.664
$ws.Range("M139").Value = -44860
$ws.Range("N139").Value = -79609.664
$ws.Range("H140").Value = 119880.664
$ws.Range("J140").Value = 119880.664
$ws.Range("L140").Value = 119880.664
$ws.Range("N140").Value = -130240.664

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 2381.8823
$ws.Range("J5").Value = 2644.3635
$ws.Range("L5").Value = 7933.0905
$ws.Range("N5").Value = -8157.0905
$ws.Range("H23").Value = 80.5
$ws.Range("I23").Value = 75
$ws.Range("J23").Value = 86
$ws.Range("K23").Value = 225
$ws.Range("L23").Value = 258
$ws.Range("M23").Value = 10
$ws.Range("N23").Value = -728
$ws.Range("H56").Value = 7000
$ws.Range("I56").Value = 7000
$ws.Range("K56").Value = 7000
$ws.Range("M56").Value = -6470
$ws.Range("H109").Value = 3144.3157
$ws.Range("J109").Value = 4399.1665
$ws.Range("L109").Value = 13197.4995
$ws.Range("N109").Value = -15277.4995
$ws.Range("H131").Value = 711665.7
$ws.Range("J131").Value = 3998.8
$ws.Range("L131").Value = 11996.4
$ws.Range("N131").Value = -22076.4
$ws.Range("H132").Value = 3138
$ws.Range("I132").Value = 3462.8
$ws.Range("J132").Value = 2867.3333
$ws.Range("K132").Value = 31165.2
$ws.Range("L132").Value = 25805.9997
$ws.Range("M132").Value = -28635.2
$ws.Range("N132").Value = -30865.9997
$ws.Range("H135").Value = 2381.8823
$ws.Range("J135").Value = 2644.3635
$ws.Range("L135").Value = 23799.2715
$ws.Range("N135").Value = -28869.2715
$ws.Range("H140").Value = 3694.6428
$ws.Range("I140").Value = 3694.6428
$ws.Range("K140").Value = 11083.9284
$ws.Range("M140").Value = -5903.928400000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 188.15384
$ws.Range("I2").Value = 195.5
$ws.Range("J2").Value = 100
$ws.Range("K2").Value = 195.5
$ws.Range("L2").Value = 100
$ws.Range("M2").Value = -82.5
$ws.Range("N2").Value = -326
$ws.Range("H10").Value = 16000000
$ws.Range("I10").Value = 16000000
$ws.Range("K10").Value = 16000000
$ws.Range("M10").Value = -15999831
$ws.Range("H15").Value = 19975
$ws.Range("J15").Value = 19975
$ws.Range("L15").Value = 19975
$ws.Range("N15").Value = -20551
$ws.Range("H18").Value = 22000
$ws.Range("J18").Value = 22000
$ws.Range("L18").Value = 22000
$ws.Range("N18").Value = -22586
$ws.Range("H19").Value = 3367334.8
$ws.Range("I19").Value = 51002
$ws.Range("K19").Value = 51002
$ws.Range("M19").Value = -50714
$ws.Range("H26").Value = 45005
$ws.Range("J26").Value = 45005
$ws.Range("L26").Value = 45005
$ws.Range("N26").Value = -45565
$ws.Range("H29").Value = 2000
$ws.Range("I29").Value = 2000
$ws.Range("K29").Value = 2000
$ws.Range("M29").Value = -1710
$ws.Range("H50").Value = 45005
$ws.Range("J50").Value = 45005
$ws.Range("L50").Value = 45005
$ws.Range("N50").Value = -46001
$ws.Range("H70").Value = 7075.778
$ws.Range("I70").Value = 6987.6665
$ws.Range("J70").Value = 7163.8887
$ws.Range("K70").Value = 6987.6665
$ws.Range("L70").Value = 7163.8887
$ws.Range("M70").Value = -6717.6665
$ws.Range("N70").Value = -7703.8887
$ws.Range("H73").Value = 7075.778
$ws.Range("I73").Value = 6987.6665
$ws.Range("J73").Value = 7163.8887
$ws.Range("K73").Value = 6987.6665
$ws.Range("L73").Value = 7163.8887
$ws.Range("M73").Value = -6051.6665
$ws.Range("N73").Value = -9035.8887
$ws.Range("H80").Value = 24199.6
$ws.Range("I80").Value = 7000
$ws.Range("J80").Value = 35666
$ws.Range("K80").Value = 7000
$ws.Range("L80").Value = 35666
$ws.Range("M80").Value = -6002
$ws.Range("N80").Value = -37662
$ws.Range("H81").Value = 19975
$ws.Range("J81").Value = 19975
$ws.Range("L81").Value = 19975
$ws.Range("N81").Value = -21971
$ws.Range("H83").Value = 24199.6
$ws.Range("I83").Value = 7000
$ws.Range("J83").Value = 35666
$ws.Range("K83").Value = 35000
$ws.Range("L83").Value = 178330
$ws.Range("M83").Value = -30008
$ws.Range("N83").Value = -188314
$ws.Range("H84").Value = 19975
$ws.Range("J84").Value = 19975
$ws.Range("L84").Value = 59925
$ws.Range("N84").Value = -69909
$ws.Range("H97").Value = 2117.05
$ws.Range("I97").Value = 2139.5
$ws.Range("J97").Value = 2027.25
$ws.Range("K97").Value = 2139.5
$ws.Range("L97").Value = 2027.25
$ws.Range("M97").Value = -1643.5
$ws.Range("N97").Value = -3019.25
$ws.Range("H98").Value = 23500
$ws.Range("J98").Value = 23500
$ws.Range("L98").Value = 23500
$ws.Range("N98").Value = -29490
$ws.Range("H112").Value = 15000
$ws.Range("J112").Value = 15000
$ws.Range("L112").Value = 15000
$ws.Range("N112").Value = -17216
$ws.Range("H113").Value = 2438
$ws.Range("I113").Value = 2516.25
$ws.Range("K113").Value = 2516.25
$ws.Range("M113").Value = -346.25
$ws.Range("H122").Value = 2413.9644
$ws.Range("I122").Value = 2107.7727
$ws.Range("K122").Value = 6323.3181
$ws.Range("M122").Value = -3873.3181
$ws.Range("H126").Value = 3299.3
$ws.Range("I126").Value = 3299.3
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 9897.900000000001
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -7427.900000000001
$ws.Range("N126").ClearContents()
$ws.Range("H132").Value = 3930.75
$ws.Range("I132").Value = 2856.7058
$ws.Range("J132").Value = 6539.143
$ws.Range("K132").Value = 8570.117400000001
$ws.Range("L132").Value = 19617.429
$ws.Range("M132").Value = -6040.117400000001
$ws.Range("N132").Value = -24677.429

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H12").Value = 41154596
$ws.Range("J12").Value = 3500
$ws.Range("L12").Value = 3500
$ws.Range("N12").Value = -3840
$ws.Range("H17").Value = 6750
$ws.Range("I17").Value = 3000
$ws.Range("J17").Value = 10500
$ws.Range("K17").Value = 3000
$ws.Range("L17").Value = 10500
$ws.Range("M17").Value = -2830
$ws.Range("N17").Value = -10840
$ws.Range("H22").Value = 1405.091
$ws.Range("J22").Value = 1715.8462
$ws.Range("L22").Value = 1715.8462
$ws.Range("N22").Value = -2305.8462
$ws.Range("H26").Value = 24999.5
$ws.Range("I26").Value = 24999.5
$ws.Range("K26").Value = 24999.5
$ws.Range("M26").Value = -24704.5
$ws.Range("H27").Value = 1405.091
$ws.Range("J27").Value = 1715.8462
$ws.Range("L27").Value = 1715.8462
$ws.Range("N27").Value = -1929.8462
$ws.Range("H46").Value = 4990.909
$ws.Range("I46").Value = 2334.3333
$ws.Range("K46").Value = 2334.3333
$ws.Range("M46").Value = -2146.3333
$ws.Range("H82").Value = 8507.286
$ws.Range("J82").Value = 7564.2856
$ws.Range("L82").Value = 7564.2856
$ws.Range("N82").Value = -8286.285599999999
$ws.Range("H85").Value = 8507.286
$ws.Range("J85").Value = 7564.2856
$ws.Range("L85").Value = 7564.2856
$ws.Range("N85").Value = -10060.2856
$ws.Range("H107").Value = 2700
$ws.Range("I107").Value = 2700
$ws.Range("K107").Value = 2700
$ws.Range("M107").Value = -780
$ws.Range("H119").Value = 105333
$ws.Range("J119").Value = 105333
$ws.Range("L119").Value = 105333
$ws.Range("N119").Value = -115009
$ws.Range("H132").Value = 2632.4443
$ws.Range("I132").Value = 2680.875
$ws.Range("J132").Value = 2245
$ws.Range("K132").Value = 8042.625
$ws.Range("L132").Value = 6735
$ws.Range("M132").Value = -5512.625
$ws.Range("N132").Value = -11795
$ws.Range("H136").Value = 6545.7144
$ws.Range("J136").Value = 3605
$ws.Range("L136").Value = 10815
$ws.Range("N136").Value = -15915

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H10").Value = 200
$ws.Range("J10").Value = 200
$ws.Range("L10").Value = 200
$ws.Range("N10").Value = -538
$ws.Range("H64").Value = 0
$ws.Range("J64").Value = 0
$ws.Range("L64").Value = 0
$ws.Range("N64").ClearContents()
$ws.Range("H67").Value = 0
$ws.Range("J67").Value = 0
$ws.Range("L67").Value = 0
$ws.Range("N67").ClearContents()
$ws.Range("H81").Value = 5047.385
$ws.Range("I81").Value = 5047.385
$ws.Range("J81").Value = 0
$ws.Range("K81").Value = 10094.77
$ws.Range("L81").Value = 0
$ws.Range("M81").Value = -9033.77
$ws.Range("N81").ClearContents()
$ws.Range("H84").Value = 5047.385
$ws.Range("I84").Value = 5047.385
$ws.Range("J84").Value = 0
$ws.Range("K84").Value = 50473.85000000001
$ws.Range("L84").Value = 0
$ws.Range("M84").Value = -45169.85000000001
$ws.Range("N84").ClearContents()
$ws.Range("H121").Value = 0
$ws.Range("J121").Value = 0
$ws.Range("L121").Value = 0
$ws.Range("N121").ClearContents()
$ws.Range("H122").Value = 106726.3
$ws.Range("I122").Value = 189071.73
$ws.Range("K122").Value = 567215.1900000001
$ws.Range("M122").Value = -564765.1900000001
$ws.Range("H126").Value = 457211.72
$ws.Range("I126").Value = 2815.5
$ws.Range("K126").Value = 8446.5
$ws.Range("M126").Value = -5976.5
$ws.Range("H132").Value = 6568.05
$ws.Range("I132").Value = 6568.05
$ws.Range("K132").Value = 19704.15
$ws.Range("M132").Value = -17174.15
$ws.Range("H135").Value = 70670.28999999999
$ws.Range("J135").Value = 70670.28999999999
$ws.Range("L135").Value = 70670.28999999999
$ws.Range("N135").Value = -80810.28999999999
$ws.Range("H136").Value = 26044.424
$ws.Range("I136").Value = 26044.424
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 78133.272
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -75583.272
$ws.Range("N136").ClearContents()
